# Applies updated win-probability matrix values (games pulled March 7)
# to team_specific_matrix/TCU_B.xlsx sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1766666666666667
$ws.Range("C2").Value = 0.5600000000000001
$ws.Range("J2").Value = 0.01333333333333333
$ws.Range("P2").Value = 0.12
$ws.Range("S2").Value = 0.13
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.04040404040404041
$ws.Range("J3").Value = 0.04040404040404041
$ws.Range("P3").Value = 0.702020202020202
$ws.Range("S3").Value = 0.202020202020202
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.6122448979591837
$ws.Range("S4").Value = 0.3265306122448979
$ws.Range("B6").Value = 0.1059907834101382
$ws.Range("D6").Value = 0.009216589861751152
$ws.Range("F6").Value = 0.08755760368663594
$ws.Range("J6").Value = 0.2396313364055299
$ws.Range("O6").Value = 0.03225806451612903
$ws.Range("Q6").Value = 0.1244239631336406
$ws.Range("R6").Value = 0.06451612903225806
$ws.Range("S6").Value = 0.336405529953917
$ws.Range("B7").Value = 0.1100478468899522
$ws.Range("D7").Value = 0.02870813397129187
$ws.Range("F7").Value = 0.07177033492822966
$ws.Range("J7").Value = 0.1052631578947368
$ws.Range("O7").Value = 0.01913875598086124
$ws.Range("Q7").Value = 0.2009569377990431
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.3732057416267943
$ws.Range("B8").Value = 0.0852994555353902
$ws.Range("D8").Value = 0.02540834845735027
$ws.Range("F8").Value = 0.07078039927404718
$ws.Range("J8").Value = 0.117967332123412
$ws.Range("O8").Value = 0.01633393829401089
$ws.Range("Q8").Value = 0.1814882032667877
$ws.Range("R8").Value = 0.05989110707803993
$ws.Range("S8").Value = 0.4428312159709619
$ws.Range("B9").Value = 0.1055555555555556
$ws.Range("D9").Value = 0.01666666666666667
$ws.Range("F9").Value = 0.04444444444444445
$ws.Range("J9").Value = 0.1388888888888889
$ws.Range("O9").Value = 0.01111111111111111
$ws.Range("Q9").Value = 0.1888888888888889
$ws.Range("R9").Value = 0.07222222222222222
$ws.Range("S9").Value = 0.4222222222222222
$ws.Range("B10").Value = 0.09822866344605476
$ws.Range("D10").Value = 0.02093397745571659
$ws.Range("E10").Value = 0.0008051529790660225
$ws.Range("F10").Value = 0.05877616747181964
$ws.Range("J10").Value = 0.1256038647342995
$ws.Range("O10").Value = 0.02012882447665056
$ws.Range("Q10").Value = 0.2230273752012882
$ws.Range("R10").Value = 0.0748792270531401
$ws.Range("S10").Value = 0.3776167471819646
$ws.Range("F11").Value = 0.003076923076923077
$ws.Range("G11").Value = 0.1384615384615385
$ws.Range("J11").Value = 0.1015384615384615
$ws.Range("K11").Value = 0.1907692307692308
$ws.Range("L11").Value = 0.5384615384615384
$ws.Range("S11").Value = 0.02769230769230769
$ws.Range("G12").Value = 0.7282608695652174
$ws.Range("J12").Value = 0.1847826086956522
$ws.Range("K12").Value = 0.02173913043478261
$ws.Range("L12").Value = 0.03804347826086957
$ws.Range("S12").Value = 0.02717391304347826
$ws.Range("G13").Value = 0.8333333333333334
$ws.Range("J13").Value = 0.1458333333333333
$ws.Range("S13").Value = 0.02083333333333333
$ws.Range("F15").Value = 0.02222222222222222
$ws.Range("H15").Value = 0.16
$ws.Range("I15").Value = 0.07555555555555556
$ws.Range("J15").Value = 0.3422222222222222
$ws.Range("K15").Value = 0.07111111111111111
$ws.Range("M15").Value = 0.004444444444444444
$ws.Range("O15").Value = 0.05777777777777778
$ws.Range("S15").Value = 0.2666666666666667
$ws.Range("F16").Value = 0.02564102564102564
$ws.Range("H16").Value = 0.1948717948717949
$ws.Range("I16").Value = 0.06153846153846154
$ws.Range("J16").Value = 0.3692307692307693
$ws.Range("K16").Value = 0.1076923076923077
$ws.Range("M16").Value = 0.005128205128205128
$ws.Range("O16").Value = 0.06666666666666667
$ws.Range("S16").Value = 0.1692307692307692
$ws.Range("F17").Value = 0.02474226804123711
$ws.Range("H17").Value = 0.1690721649484536
$ws.Range("I17").Value = 0.1175257731958763
$ws.Range("J17").Value = 0.4020618556701031
$ws.Range("K17").Value = 0.08865979381443299
$ws.Range("M17").Value = 0.02061855670103093
$ws.Range("O17").Value = 0.05360824742268041
$ws.Range("S17").Value = 0.1237113402061856
$ws.Range("F18").Value = 0.01136363636363636
$ws.Range("H18").Value = 0.1761363636363636
$ws.Range("I18").Value = 0.07386363636363637
$ws.Range("J18").Value = 0.3636363636363636
$ws.Range("K18").Value = 0.1079545454545455
$ws.Range("M18").Value = 0.01704545454545454
$ws.Range("O18").Value = 0.05681818181818182
$ws.Range("S18").Value = 0.1931818181818182
$ws.Range("F19").Value = 0.01583873290136789
$ws.Range("H19").Value = 0.2613390928725702
$ws.Range("I19").Value = 0.05759539236861051
$ws.Range("J19").Value = 0.3383729301655867
$ws.Range("K19").Value = 0.1159107271418287
$ws.Range("M19").Value = 0.02375809935205184
$ws.Range("N19").Value = 0.001439884809215263
$ws.Range("O19").Value = 0.05903527717782577
$ws.Range("S19").Value = 0.1267098632109431
